$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $newVal) {
    $rng = $ws.Range($cellRef)
    if ($newVal -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # New value reads as a plain number (e.g. "1.00", "0.161").
        # Force text storage first so Excel keeps the literal digits
        # (incl. trailing zeros / multi-dot thousands groups) instead of
        # silently coercing to a numeric cell, then drop the number-format
        # bump back to the default style so only the value changes.
        $rng.NumberFormat = "@"
        $rng.Value = $newVal
        $rng.Style = "Normal"
    } else {
        $rng.Value = $newVal
    }
}

# (cell reference, new text value) pairs taken from the source update.
# Each pair is prefixed with the unary "," so the outer array stays an
# array-of-arrays instead of being flattened into one long flat list.
$updates = @(
    ,@('D2', '68.063.39')
    ,@('E2', '  -1.58%  ')
    ,@('D3', '2.412.53')
    ,@('E3', '  -2.59%  ')
    ,@('E4', '  -0.02%  ')
    ,@('D5', '554.84')
    ,@('E5', '  -1.47%  ')
    ,@('D6', '159.09')
    ,@('E6', '  -2.48%  ')
    ,@('E7', '  -0.03%  ')
    ,@('E8', '  -0.60%  ')
    ,@('D9', '0.161')
    ,@('E9', '  +6.10%  ')
    ,@('E10', '  -1.79%  ')
    ,@('E11', '  -1.27%  ')
    ,@('D12', '4.63')
    ,@('E12', '  -5.06%  ')
    ,@('D13', '67.947.82')
    ,@('D14', '2.854.85')
    ,@('E14', '  -1.75%  ')
    ,@('E15', '  +1.61%  ')
    ,@('D16', '22.77')
    ,@('E16', '  -4.08%  ')
    ,@('D17', '2.413.74')
    ,@('E17', '  -3.51%  ')
    ,@('D18', '10.39')
    ,@('E18', '  -3.87%  ')
    ,@('D19', '331.21')
    ,@('E19', '  -2.55%  ')
    ,@('E20', '  -2.78%  ')
    ,@('B22', 'Dai')
    ,@('C22', 'https://coinranking.com/coin/MoTuySvg7+dai-dai')
    ,@('D22', '1.00')
    ,@('E22', '  +0.00%  ')
    ,@('B23', 'SuiNetwork')
    ,@('C23', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui')
    ,@('D23', '1.88')
    ,@('E23', '  -1.52%  ')
    ,@('D24', '66.19')
    ,@('E24', '  -1.73%  ')
    ,@('E25', '  -1.68%  ')
    ,@('D26', '2.537.80')
    ,@('E26', '  -2.56%  ')
    ,@('E27', '  -1.42%  ')
    ,@('E28', '  -2.18%  ')
    ,@('E29', '  -1.87%  ')
    ,@('D30', '1.00')
    ,@('E30', '  +0.04%  ')
    ,@('D31', '421.65')
    ,@('E31', '  -3.16%  ')
    ,@('E32', '  -1.59%  ')
    ,@('E33', '  -2.00%  ')
    ,@('D34', '159.10')
    ,@('E34', '  +0.94%  ')
    ,@('D35', '19.02')
    ,@('E35', '  -0.24%  ')
    ,@('E36', '  -0.01%  ')
    ,@('E37', '  -0.40%  ')
    ,@('D38', '0.104')
    ,@('E38', '  -5.48%  ')
    ,@('E39', '  -2.22%  ')
    ,@('E40', '  -3.70%  ')
    ,@('E41', '  -0.95%  ')
    ,@('E42', '  -1.55%  ')
    ,@('D43', '132.36')
    ,@('E43', '  -1.04%  ')
    ,@('E44', '  -1.60%  ')
    ,@('E45', '  -5.43%  ')
    ,@('D46', '0.0712')
    ,@('E46', '  -0.86%  ')
    ,@('E47', '  -1.79%  ')
    ,@('D48', '0.553')
    ,@('E48', '  -2.15%  ')
    ,@('E49', '  -0.37%  ')
    ,@('E50', '  -1.11%  ')
    ,@('E51', '  -3.53%  ')
)

foreach ($u in $updates) {
    Set-CellText $u[0] $u[1]
}
